$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '93.192.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.434.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '620.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.38'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.393'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.75%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.430.82'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.26%  '
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.072.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '93.059.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.431.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '500.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.445'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("E26").Value = '  -4.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '92.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.618.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.137'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.172'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '30.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '554.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.19%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.40'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.921'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.149'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.13%  '
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("E45").Value = '  -1.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0410'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("E50").Value = '  -4.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.38%  '
